$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 37040
$ws.Range("E2").Value = 724605718990
$ws.Range("F2").Value = 12522762620
$ws.Range("G2").Value = -0.16848

$ws.Range("D3").Value = 2041.15
$ws.Range("E3").Value = 246157223708
$ws.Range("F3").Value = 11951582315
$ws.Range("G3").Value = -1.03607

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 86823534820
$ws.Range("F4").Value = 34449847337
$ws.Range("G4").Value = -0.09526

$ws.Range("D5").Value = 248.65
$ws.Range("E5").Value = 38317048614
$ws.Range("F5").Value = 879651128
$ws.Range("G5").Value = -0.40342

$ws.Range("D6").Value = 0.664659
$ws.Range("E6").Value = 35757221274
$ws.Range("F6").Value = 1365124642
$ws.Range("G6").Value = -1.08723

$ws.Range("B7").Value = "SOL"
$ws.Range("C7").Value = "Solana"
$ws.Range("D7").Value = 58.68
$ws.Range("E7").Value = 25117403863
$ws.Range("F7").Value = 4628222328
$ws.Range("G7").Value = 4.41822

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "USDC"
$ws.Range("D8").Value = 0.999221
$ws.Range("E8").Value = 24224936652
$ws.Range("F8").Value = 6413967726
$ws.Range("G8").Value = -0.16951

$ws.Range("D9").Value = 2045.19
$ws.Range("E9").Value = 18343678764
$ws.Range("F9").Value = 12354607
$ws.Range("G9").Value = -0.79071

$ws.Range("D10").Value = 0.384794
$ws.Range("E10").Value = 13504796837
$ws.Range("F10").Value = 538761565
$ws.Range("G10").Value = 1.52594

$ws.Range("D11").Value = 0.078669
$ws.Range("E11").Value = 11209518930
$ws.Range("F11").Value = 1452450603
$ws.Range("G11").Value = -1.92124

$ws.Range("B12").Value = "TRX"
$ws.Range("C12").Value = "TRON"
$ws.Range("D12").Value = 0.108046
$ws.Range("E12").Value = 9592159930
$ws.Range("F12").Value = 296989630
$ws.Range("G12").Value = 1.23702

$ws.Range("B13").Value = "LINK"
$ws.Range("C13").Value = "Chainlink"
$ws.Range("D13").Value = 15.91
$ws.Range("E13").Value = 8921358550
$ws.Range("F13").Value = 1392273390
$ws.Range("G13").Value = 4.70082

$ws.Range("B14").Value = "MATIC"
$ws.Range("C14").Value = "Polygon"
$ws.Range("D14").Value = 0.832392
$ws.Range("E14").Value = 7730131211
$ws.Range("F14").Value = 696743261
$ws.Range("G14").Value = 2.1093

$ws.Range("B15").Value = "DOT"
$ws.Range("C15").Value = "Polkadot"
$ws.Range("D15").Value = 5.77
$ws.Range("E15").Value = 7507447148
$ws.Range("F15").Value = 536924480
$ws.Range("G15").Value = 8.306940000000001

$ws.Range("B16").Value = "AVAX"
$ws.Range("C16").Value = "Avalanche"
$ws.Range("D16").Value = 18.44
$ws.Range("E16").Value = 6577465473
$ws.Range("F16").Value = 1402997291
$ws.Range("G16").Value = 29.19223

$ws.Range("D17").Value = 37045
$ws.Range("E17").Value = 6068339256
$ws.Range("F17").Value = 229323005
$ws.Range("G17").Value = -0.01298

$ws.Range("D18").Value = 75.13
$ws.Range("E18").Value = 5589109160
$ws.Range("F18").Value = 604310194
$ws.Range("G18").Value = 0.82397

$ws.Range("B19").Value = "DAI"
$ws.Range("C19").Value = "Dai"
$ws.Range("D19").Value = 0.998789
$ws.Range("E19").Value = 5344232878
$ws.Range("F19").Value = 149291487
$ws.Range("G19").Value = -0.09699000000000001

$ws.Range("B20").Value = "SHIB"
$ws.Range("C20").Value = "Shiba Inu"
$ws.Range("D20").Value = 0.000008970000000000001
$ws.Range("E20").Value = 5317565933
$ws.Range("F20").Value = 316725240
$ws.Range("G20").Value = -2.31299

$ws.Range("B21").Value = "TON"
$ws.Range("C21").Value = "Toncoin"
$ws.Range("D21").Value = 2.42
$ws.Range("E21").Value = 5176194300
$ws.Range("F21").Value = 54764374
$ws.Range("G21").Value = -0.30984

$ws.Range("B22").Value = "BCH"
$ws.Range("C22").Value = "Bitcoin Cash"
$ws.Range("D22").Value = 236.81
$ws.Range("E22").Value = 4645632726
$ws.Range("F22").Value = 143715572
$ws.Range("G22").Value = -0.36195

$ws.Range("B23").Value = "UNI"
$ws.Range("C23").Value = "Uniswap"
$ws.Range("D23").Value = 5.35
$ws.Range("E23").Value = 4054865364
$ws.Range("F23").Value = 168902663
$ws.Range("G23").Value = -0.83385

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "LEO Token"
$ws.Range("D24").Value = 4.09
$ws.Range("E24").Value = 3797066618
$ws.Range("F24").Value = 1027107
$ws.Range("G24").Value = -0.25779

$ws.Range("B25").Value = "OKB"
$ws.Range("C25").Value = "OKB"
$ws.Range("D25").Value = 60.22
$ws.Range("E25").Value = 3619693066
$ws.Range("F25").Value = 14205455
$ws.Range("G25").Value = -1.30488

$ws.Range("B26").Value = "XLM"
$ws.Range("C26").Value = "Stellar"
$ws.Range("D26").Value = 0.125533
$ws.Range("E26").Value = 3520446780
$ws.Range("F26").Value = 117958647
$ws.Range("G26").Value = 0.50259

$ws.Range("B27").Value = "TUSD"
$ws.Range("C27").Value = "TrueUSD"
$ws.Range("D27").Value = 0.997806
$ws.Range("E27").Value = 3334473511
$ws.Range("F27").Value = 137459642
$ws.Range("G27").Value = -0.24576

$ws.Range("B28").Value = "XMR"
$ws.Range("C28").Value = "Monero"
$ws.Range("D28").Value = 169.29
$ws.Range("E28").Value = 3078729747
$ws.Range("F28").Value = 100387755
$ws.Range("G28").Value = -1.39343

$ws.Range("D29").Value = 20.08
$ws.Range("E29").Value = 2884038090
$ws.Range("F29").Value = 256770497
$ws.Range("G29").Value = 0.01911

$ws.Range("B30").Value = "CRO"
$ws.Range("C30").Value = "Cronos"
$ws.Range("D30").Value = 0.107997
$ws.Range("E30").Value = 2870973380
$ws.Range("F30").Value = 98245056
$ws.Range("G30").Value = 2.31098

$ws.Range("B31").Value = "ATOM"
$ws.Range("C31").Value = "Cosmos Hub"
$ws.Range("D31").Value = 9.41
$ws.Range("E31").Value = 2770617037
$ws.Range("F31").Value = 187429296
$ws.Range("G31").Value = 3.5936

$ws.Range("B32").Value = "FIL"
$ws.Range("C32").Value = "Filecoin"
$ws.Range("D32").Value = 4.79
$ws.Range("E32").Value = 2249195580
$ws.Range("F32").Value = 265052488
$ws.Range("G32").Value = 4.02152

$ws.Range("D33").Value = 0.062824
$ws.Range("E33").Value = 2116099234
$ws.Range("F33").Value = 91131862
$ws.Range("G33").Value = 0.03713

$ws.Range("D34").Value = 4.53
$ws.Range("E34").Value = 2044221073
$ws.Range("F34").Value = 69334403
$ws.Range("G34").Value = 3.18482

$ws.Range("B35").Value = "APT"
$ws.Range("C35").Value = "Aptos"
$ws.Range("D35").Value = 8.130000000000001
$ws.Range("E35").Value = 2040217993
$ws.Range("F35").Value = 320428875
$ws.Range("G35").Value = 10.08906

$ws.Range("D36").Value = 2.22
$ws.Range("E36").Value = 1987240220
$ws.Range("F36").Value = 113720477
$ws.Range("G36").Value = -2.3128

$ws.Range("B37").Value = "KAS"
$ws.Range("C37").Value = "Kaspa"
$ws.Range("D37").Value = 0.08903899999999999
$ws.Range("E37").Value = 1921720603
$ws.Range("F37").Value = 40047333
$ws.Range("G37").Value = 2.48551

$ws.Range("B38").Value = "BUSD"
$ws.Range("C38").Value = "BUSD"
$ws.Range("D38").Value = 0.99847
$ws.Range("E38").Value = 1876120940
$ws.Range("F38").Value = 1146168399
$ws.Range("G38").Value = -0.238

$ws.Range("B39").Value = "NEAR"
$ws.Range("C39").Value = "NEAR Protocol"
$ws.Range("D39").Value = 1.77
$ws.Range("E39").Value = 1768201614
$ws.Range("F39").Value = 325462030
$ws.Range("G39").Value = 13.83276

$ws.Range("B40").Value = "VET"
$ws.Range("C40").Value = "VeChain"
$ws.Range("D40").Value = 0.02230666
$ws.Range("E40").Value = 1627184770
$ws.Range("F40").Value = 56816636
$ws.Range("G40").Value = -0.80003

$ws.Range("B41").Value = "RUNE"
$ws.Range("C41").Value = "THORChain"
$ws.Range("D41").Value = 5.29
$ws.Range("E41").Value = 1613622408
$ws.Range("F41").Value = 1329471861
$ws.Range("G41").Value = 20.84753

$ws.Range("B42").Value = "QNT"
$ws.Range("C42").Value = "Quant"
$ws.Range("D42").Value = 106.25
$ws.Range("E42").Value = 1551057007
$ws.Range("F42").Value = 31319853
$ws.Range("G42").Value = -4.27393

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "OP"
$ws.Range("C43").Value = "Optimism"
$ws.Range("D43").Value = 1.71
$ws.Range("E43").Value = 1513960032
$ws.Range("F43").Value = 166336923
$ws.Range("G43").Value = 1.24694

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "INJ"
$ws.Range("C44").Value = "Injective"
$ws.Range("D44").Value = 17.61
$ws.Range("E44").Value = 1481359860
$ws.Range("F44").Value = 148087340
$ws.Range("G44").Value = -3.02037

$ws.Range("B45").Value = "ARB"
$ws.Range("C45").Value = "Arbitrum"
$ws.Range("D45").Value = 1.14
$ws.Range("E45").Value = 1454452587
$ws.Range("F45").Value = 418265046
$ws.Range("G45").Value = -1.01507

$ws.Range("B46").Value = "IMX"
$ws.Range("C46").Value = "Immutable"
$ws.Range("D46").Value = 1.13
$ws.Range("E46").Value = 1434707009
$ws.Range("F46").Value = 785432074
$ws.Range("G46").Value = 6.84381

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "AAVE"
$ws.Range("C47").Value = "Aave"
$ws.Range("D47").Value = 96.59999999999999
$ws.Range("E47").Value = 1421197962
$ws.Range("F47").Value = 202826884
$ws.Range("G47").Value = -0.37999

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "MNT"
$ws.Range("C48").Value = "Mantle"
$ws.Range("D48").Value = 0.438165
$ws.Range("E48").Value = 1362766826
$ws.Range("F48").Value = 7078277
$ws.Range("G48").Value = -3.38986

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "GRT"
$ws.Range("C49").Value = "The Graph"
$ws.Range("D49").Value = 0.137074
$ws.Range("E49").Value = 1279826091
$ws.Range("F49").Value = 123752311
$ws.Range("G49").Value = 2.20593

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "RETH"
$ws.Range("C50").Value = "Rocket Pool ETH"
$ws.Range("D50").Value = 2233.25
$ws.Range("E50").Value = 1197182243
$ws.Range("F50").Value = 7983732
$ws.Range("G50").Value = -0.61734

$ws.Range("B51").Value = "MKR"
$ws.Range("C51").Value = "Maker"
$ws.Range("D51").Value = 1293.12
$ws.Range("E51").Value = 1168378184
$ws.Range("F51").Value = 86633999
$ws.Range("G51").Value = -0.483
